$d = $word.ActiveDocument

# Locate the "References" heading paragraph.
$refPara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    $trimmed = $txt.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "References") {
        $refPara = $p
        break
    }
}

if ($refPara -ne $null) {
    # Remove everything from the start of the "References" heading through
    # the end of the document content (heading + citation + trailing blank
    # paragraph), leaving the preceding blank-space paragraph intact.
    $startPos = $refPara.Range.Start
    $endPos = $d.Content.End
    $killRange = $d.Range($startPos, $endPos)
    $killRange.Delete()
}
